$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("ZoneLetter") holds "V" for data rows 2-180.
# Change every "V" to "T" (subdivision of area 2 according to time class).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 180 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "V") {
        $cell.Value2 = "T"
    }
}
